$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 17
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = ""  # was -388
$ws.Range("H74").Value = 3053.818
$ws.Range("I74").Value = 3179.2
$ws.Range("K74").Value = 3179.2
$ws.Range("M74").Value = -2243.2
$ws.Range("H77").Value = 3053.818
$ws.Range("I77").Value = 3179.2
$ws.Range("K77").Value = 15896
$ws.Range("M77").Value = -11216
$ws.Range("H99").Value = 607.1429000000001
$ws.Range("J99").Value = 1111
$ws.Range("L99").Value = 3333
$ws.Range("N99").Value = -6329
$ws.Range("H100").Value = 3578.5557
$ws.Range("I100").Value = 2721.3333
$ws.Range("J100").Value = 4007.1667
$ws.Range("K100").Value = 2721.3333
$ws.Range("L100").Value = 4007.1667
$ws.Range("M100").Value = -2180.3333
$ws.Range("N100").Value = -5089.1667
$ws.Range("H132").Value = 1553.75
$ws.Range("I132").Value = 1553.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4661.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""  # was -1202
$ws.Range("N132").Value = -2131.25
$ws.Range("H138").Value = 1451237.1
$ws.Range("I138").Value = 885.5
$ws.Range("J138").Value = 3033439
$ws.Range("K138").Value = 2656.5
$ws.Range("L138").Value = 9100317
$ws.Range("M138").Value = 2483.5
$ws.Range("N138").Value = -9110597

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 38544424
$ws.Range("I61").Value = 83335570
$ws.Range("K61").Value = 83335570
$ws.Range("M61").Value = -83335358
$ws.Range("H74").Value = 6256513.5
$ws.Range("I74").Value = 9260844
$ws.Range("J74").Value = 16750.309
$ws.Range("K74").Value = 9260844
$ws.Range("L74").Value = 16750.309
$ws.Range("M74").Value = -9259970
$ws.Range("N74").Value = -18498.309
$ws.Range("H77").Value = 6256513.5
$ws.Range("I77").Value = 9260844
$ws.Range("J77").Value = 16750.309
$ws.Range("K77").Value = 46304220
$ws.Range("L77").Value = 83751.54500000001
$ws.Range("M77").Value = -46299852
$ws.Range("N77").Value = -92487.54500000001
$ws.Range("H88").Value = 1713.1428
$ws.Range("I88").Value = 1555.5714
$ws.Range("J88").Value = 1870.7142
$ws.Range("K88").Value = 1555.5714
$ws.Range("L88").Value = 1870.7142
$ws.Range("M88").Value = -1149.5714
$ws.Range("N88").Value = -2682.7142
$ws.Range("H91").Value = 1713.1428
$ws.Range("I91").Value = 1555.5714
$ws.Range("J91").Value = 1870.7142
$ws.Range("K91").Value = 1555.5714
$ws.Range("L91").Value = 1870.7142
$ws.Range("M91").Value = -151.5714
$ws.Range("N91").Value = -4678.7142
$ws.Range("H102").Value = 11146.333
$ws.Range("I102").Value = 13745.182
$ws.Range("K102").Value = 13745.182
$ws.Range("M102").Value = -12123.182
$ws.Range("H110").Value = 5017.3335
$ws.Range("I110").Value = 5032
$ws.Range("J110").Value = 4900
$ws.Range("K110").Value = 5032
$ws.Range("L110").Value = 4900
$ws.Range("M110").Value = -2987
$ws.Range("N110").Value = -8990
$ws.Range("H136").Value = 38544424
$ws.Range("I136").Value = 83335570
$ws.Range("K136").Value = 250006710
$ws.Range("M136").Value = -250004160

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 25000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 25000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = ""  # was 0
$ws.Range("M31").Value = 25000
$ws.Range("N31").Value = -25504
$ws.Range("H75").Value = 99995
$ws.Range("J75").Value = 99995
$ws.Range("L75").Value = 99995
$ws.Range("N75").Value = -101867
$ws.Range("H78").Value = 99995
$ws.Range("J78").Value = 99995
$ws.Range("L78").Value = 299985
$ws.Range("N78").Value = -309345
$ws.Range("H94").Value = 2026.5294
$ws.Range("I94").Value = 1788
$ws.Range("K94").Value = 1788
$ws.Range("M94").Value = -1337
$ws.Range("H105").Value = 1916.1666
$ws.Range("I105").Value = 1799.4
$ws.Range("K105").Value = 1799.4
$ws.Range("M105").Value = -52.40000000000009
$ws.Range("H107").Value = 1804.75
$ws.Range("I107").Value = 1487.5
$ws.Range("K107").Value = 1487.5
$ws.Range("M107").Value = 432.5
$ws.Range("H134").Value = 49312.793
$ws.Range("I134").Value = 6349.7
$ws.Range("K134").Value = 19049.1
$ws.Range("M134").Value = -16514.1

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 35134.332
$ws.Range("I93").Value = 35134.332
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 35134.332
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = ""  # was -23127
$ws.Range("N93").Value = -33262.332
$ws.Range("H105").Value = 2544
$ws.Range("I105").Value = 2161.8
$ws.Range("K105").Value = 2161.8
$ws.Range("M105").Value = -414.8000000000002
$ws.Range("H125").Value = 30080.25
$ws.Range("J125").Value = 30080.25
$ws.Range("L125").Value = 30080.25
$ws.Range("N125").Value = -35000.25
$ws.Range("H132").Value = 2149.5715
$ws.Range("I132").Value = 1930.3077
$ws.Range("K132").Value = 5790.9231
$ws.Range("M132").Value = -3260.9231
$ws.Range("H134").Value = 1254228.2
$ws.Range("I134").Value = 1668135.9
$ws.Range("J134").Value = 12505.5
$ws.Range("K134").Value = 5004407.699999999
$ws.Range("L134").Value = 37516.5
$ws.Range("M134").Value = -5001872.699999999
$ws.Range("N134").Value = -42586.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 500002500
$ws.Range("J25").Value = 500002500
$ws.Range("L25").Value = 1500007500
$ws.Range("N25").Value = -1500007838
$ws.Range("H30").Value = 500002500
$ws.Range("J30").Value = 500002500
$ws.Range("L30").Value = 1500007500
$ws.Range("N30").Value = -1500007704
$ws.Range("H39").Value = 141054.2
$ws.Range("J39").Value = 155649.8
$ws.Range("L39").Value = 466949.4
$ws.Range("N39").Value = -467537.4
$ws.Range("H46").Value = 673.0769
$ws.Range("I46").Value = 768.625
$ws.Range("J46").Value = 520.2
$ws.Range("K46").Value = 2305.875
$ws.Range("L46").Value = 1560.6
$ws.Range("M46").Value = -2214.875
$ws.Range("N46").Value = -1742.6
$ws.Range("H107").Value = 1058.6
$ws.Range("J107").Value = 1200.75
$ws.Range("L107").Value = 3602.25
$ws.Range("N107").Value = -7442.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 84316
$ws.Range("I82").Value = 42300
$ws.Range("J82").Value = 105324
$ws.Range("K82").Value = 42300
$ws.Range("L82").Value = 105324
$ws.Range("M82").Value = -41917
$ws.Range("N82").Value = -106090
$ws.Range("H85").Value = 84316
$ws.Range("I85").Value = 42300
$ws.Range("J85").Value = 105324
$ws.Range("K85").Value = 42300
$ws.Range("L85").Value = 105324
$ws.Range("M85").Value = -40974
$ws.Range("N85").Value = -107976
$ws.Range("H122").Value = 3569
$ws.Range("I122").Value = 2997.1667
$ws.Range("K122").Value = 8991.500100000001
$ws.Range("M122").Value = -6541.500100000001
$ws.Range("H126").Value = 10071
$ws.Range("I126").Value = 7624.25
$ws.Range("J126").Value = 13333.333
$ws.Range("K126").Value = 22872.75
$ws.Range("L126").Value = 39999.999
$ws.Range("M126").Value = -20402.75
$ws.Range("N126").Value = -44939.999
$ws.Range("H132").Value = 333338340
$ws.Range("I132").Value = 500005000
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 1500015000
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1500012470
$ws.Range("N132").Value = -20060

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1456.4375
$ws.Range("I16").Value = 1650.9231
$ws.Range("K16").Value = 1650.9231
$ws.Range("M16").Value = -1480.9231
$ws.Range("H132").Value = 369611.1
$ws.Range("I132").Value = 716964
$ws.Range("J132").Value = 65677.31
$ws.Range("K132").Value = 2150892
$ws.Range("L132").Value = 197031.93
$ws.Range("M132").Value = -2148362
$ws.Range("N132").Value = -202091.93
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 9417052
$ws.Range("J75").Value = 9417052
$ws.Range("L75").Value = 9417052
$ws.Range("N75").Value = -9418924
$ws.Range("H78").Value = 9417052
$ws.Range("J78").Value = 9417052
$ws.Range("L78").Value = 28251156
$ws.Range("N78").Value = -28260516
$ws.Range("H81").Value = 67900
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 100850
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 201700
$ws.Range("M81").Value = -2939
$ws.Range("N81").Value = -203822
$ws.Range("H84").Value = 67900
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 100850
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 1008500
$ws.Range("M84").Value = -14696
$ws.Range("N84").Value = -1019108
$ws.Range("H86").Value = 57997.5
$ws.Range("J86").Value = 57997.5
$ws.Range("L86").Value = 57997.5
$ws.Range("N86").Value = -60243.5
$ws.Range("H89").Value = 57997.5
$ws.Range("J89").Value = 57997.5
$ws.Range("L89").Value = 289987.5
$ws.Range("N89").Value = -301219.5
$ws.Range("H141").Value = 67333
$ws.Range("J141").Value = 67333
$ws.Range("L141").Value = 67333
$ws.Range("N141").Value = -77693
